$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report week / volume number) ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Cells that must become literal TEXT values (numeric -> text "0") ---
# Copy full cell (value+style) from C14, which already holds text "0" with the right style,
# to avoid Excel auto-converting the numeric-looking string back into a number.
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("C23"))

# --- Cells that must become numeric values (previously literal TEXT placeholders) ---
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F31").Value = 1
$ws.Range("F31").NumberFormat = '#,##0'
$ws.Range("I31").Value = 1
$ws.Range("I31").NumberFormat = '#,##0'

# --- Plain numeric value updates ---
$ws.Range("L15").Value = -33.333333333333
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 133.333333333333
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = 53.333333333333
$ws.Range("L16").Value = 43.75
$ws.Range("M16").Value = -8
$ws.Range("N16").Value = -78.703703703703
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 65
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 91.176470588235
$ws.Range("L17").Value = 71.052631578947
$ws.Range("M17").Value = 109.677419354839
$ws.Range("N17").Value = -16.666666666666
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 21
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = 50
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = -89.5
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 8
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 54
$ws.Range("J19").Value = 66
$ws.Range("K19").Value = -18.181818181818
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = -11.475409836065
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 300
$ws.Range("L20").Value = 41.176470588235
$ws.Range("M20").Value = 242.857142857143
$ws.Range("N20").Value = -83.098591549295
$ws.Range("C21").Value = 9
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 65.384615384615
$ws.Range("I21").Value = 189
$ws.Range("J21").Value = 141
$ws.Range("K21").Value = 34.042553191489
$ws.Range("L21").Value = 26.845637583892
$ws.Range("M21").Value = 53.658536585365
$ws.Range("N21").Value = -68.394648829431
$ws.Range("L23").Value = 9.090909090909
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 26
$ws.Range("G24").Value = 25
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 159
$ws.Range("J24").Value = 179
$ws.Range("K24").Value = -11.173184357541
$ws.Range("L24").Value = -17.1875
$ws.Range("M24").Value = 63.917525773195
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 66.666666666666
$ws.Range("I25").Value = 73
$ws.Range("J25").Value = 112
$ws.Range("K25").Value = -34.821428571428
$ws.Range("L25").Value = -20.652173913043
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 35.294117647058
$ws.Range("I26").Value = 97
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 21.25
$ws.Range("L26").Value = 36.619718309859
$ws.Range("M26").Value = -29.710144927536
$ws.Range("L27").Value = 0
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = -20
$ws.Range("H31").Value = 0
$ws.Range("K31").Value = -50
